$d = $word.ActiveDocument

# Locate the end of the first paragraph's text ("This is a Microsoft word document.")
$para = $d.Paragraphs.Item(1)
$r = $para.Range
# Move the end back by 1 to exclude the paragraph mark, so we insert right after the period.
$r.End = $r.End - 1
$r.Collapse(0)  # wdCollapseEnd = 0

$r.InsertAfter(" (")
$r.Collapse(0)

$r.InsertAfter("Changed main")
$r.Collapse(0)

$r.InsertAfter(")")
